$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.187.43"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.768.42"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.71"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5232"
$ws.Range("E7").Value = "  +10.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3608"
$ws.Range("E8").Value = "  +5.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.51"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07337"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.080"
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9986"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.49"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.054"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").Value = "1.767.33"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.953"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.25"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06410"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9987"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.65"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.827"
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("D23").Value = "27.284.18"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.062"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.92"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.09"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "1.970.44"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.320"
$ws.Range("E29").Value = "  +11.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.25"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.059"
$ws.Range("E31").Value = "  +4.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09721"
$ws.Range("E32").Value = "  +6.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.524"
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.602"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02220"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05965"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.17"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6112"
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.806"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.431"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.965"
$ws.Range("E42").Value = "  +6.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.133"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.24"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5735"
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.617"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.92"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.874"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.106"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06694"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.74"
$ws.Range("E51").Value = "  +1.40%  "
